$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- New row 17: Athira John ---
$ws.Range("A17").Value = "Athira John"
$ws.Range("B17").Value = 36
$ws.Range("C17").Value = 38
$ws.Range("E17").Value = "Taster"

# --- New row 18: Athira Niar ---
$ws.Range("A18").Value = "Athira Niar"
$ws.Range("B18").Value = 31
$ws.Range("C18").Value = 33
$ws.Range("E18").Value = "Taster"

# --- New E14 value (SuperTaster) ---
$ws.Range("E14").Value = "SuperTaster"

# --- Rebuild the Mean column as one shared formula across D2:D18 ---
$ws.Range("D2:D18").Formula = "=AVERAGE(B2:C2)"

# --- Column E width (matches Excel's auto best-fit width for "Non Taster") ---
$ws.Columns.Item(5).ColumnWidth = 10.7369791666667

# --- Update active selection to A19, matching the post-edit cursor position ---
$ws.Range("A19").Select() | Out-Null
